$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Chillán - Zanahoria.
# It belongs chronologically right after the existing row 229, so insert a
# fresh row at 230 (this pushes the former rows 230-245 down to 231-246,
# keeping their data intact) and fill it in with the new observation.
$ws.Rows.Item(230).Insert()

$ws.Cells.Item(230, 1).Value = 7
$ws.Cells.Item(230, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(230, 3).Value = "Ñuble"
$ws.Cells.Item(230, 4).Value = 44585
$ws.Cells.Item(230, 5).Value = 16
$ws.Cells.Item(230, 6).Value = 100114013
$ws.Cells.Item(230, 7).Value = "Zanahoria"
$ws.Cells.Item(230, 8).Value = "Sin especificar"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 200
$ws.Cells.Item(230, 11).Value = 6500
$ws.Cells.Item(230, 12).Value = 7000
$ws.Cells.Item(230, 13).Value = 6750
$ws.Cells.Item(230, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(230, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(230, 16).Value = 338
$ws.Cells.Item(230, 17).Value = 20
$ws.Cells.Item(230, 18).Value = "Hortaliza"
